$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values from the repulled data source
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = -2
$ws.Range("F6").Value = -4
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = -4
$ws.Range("F9").Value = 5
$ws.Range("F11").Value = -2
$ws.Range("F12").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = -1
$ws.Range("F17").Value = -4
$ws.Range("F18").Value = -4
$ws.Range("F19").Value = -2
$ws.Range("F20").Value = 6
$ws.Range("F22").Value = 4
$ws.Range("F23").Value = 7
$ws.Range("F25").Value = -4
$ws.Range("F26").Value = -5
$ws.Range("F27").Value = -5
$ws.Range("F28").Value = -4
$ws.Range("F29").Value = 4
$ws.Range("F30").Value = 1
